$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in the header cell F1
$ws.Range("F1").Value = "Last status check on: 08.02.2022 09:00"

# Convert D8 from a text "+0.4" into a real number 0.4
$ws.Range("D8").Value = 0.4

# Convert E8 from a text date-time string into a real Excel date/time serial
# value, matching the number format used by the other rows in column E.
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 44600.36524305555
